$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Paragraphs.Item(1).Range.Text = "2023-02-17 Friday"

# Update each arithmetic-problem cell in the table directly by position
# (direct Range.Text assignment avoids any cross-cell substring collisions)
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "46-13="
$tbl.Cell(1, 2).Range.Text = "51-24="
$tbl.Cell(1, 3).Range.Text = "14+63="
$tbl.Cell(1, 4).Range.Text = "39+18="
$tbl.Cell(1, 5).Range.Text = "43-11="
$tbl.Cell(2, 1).Range.Text = "83-69="
$tbl.Cell(2, 2).Range.Text = "27-12="
$tbl.Cell(2, 3).Range.Text = "88-40="
$tbl.Cell(2, 4).Range.Text = "58+14="
$tbl.Cell(2, 5).Range.Text = "8+63="
$tbl.Cell(3, 1).Range.Text = "24-19="
$tbl.Cell(3, 2).Range.Text = "49+3="
$tbl.Cell(3, 3).Range.Text = "44-5="
$tbl.Cell(3, 4).Range.Text = "74-60="
$tbl.Cell(3, 5).Range.Text = "36-8="
$tbl.Cell(4, 1).Range.Text = "9+41="
$tbl.Cell(4, 2).Range.Text = "43-23="
$tbl.Cell(4, 3).Range.Text = "61+20="
$tbl.Cell(4, 4).Range.Text = "99-89="
$tbl.Cell(4, 5).Range.Text = "61-7="
$tbl.Cell(5, 1).Range.Text = "3+27="
$tbl.Cell(5, 2).Range.Text = "12+21="
$tbl.Cell(5, 3).Range.Text = "79+17="
$tbl.Cell(5, 4).Range.Text = "23+24="
$tbl.Cell(5, 5).Range.Text = "58-2="
$tbl.Cell(6, 1).Range.Text = "12+35="
$tbl.Cell(6, 2).Range.Text = "80-40="
$tbl.Cell(6, 3).Range.Text = "93+5="
$tbl.Cell(6, 4).Range.Text = "97-58="
$tbl.Cell(6, 5).Range.Text = "91-56="
$tbl.Cell(7, 1).Range.Text = "3+36="
$tbl.Cell(7, 2).Range.Text = "34+28="
$tbl.Cell(7, 3).Range.Text = "59+35="
$tbl.Cell(7, 4).Range.Text = "65-47="
$tbl.Cell(7, 5).Range.Text = "50-35="
$tbl.Cell(8, 1).Range.Text = "85-62="
$tbl.Cell(8, 2).Range.Text = "1+50="
$tbl.Cell(8, 3).Range.Text = "81-18="
$tbl.Cell(8, 4).Range.Text = "51-6="
$tbl.Cell(8, 5).Range.Text = "94-75="
$tbl.Cell(9, 1).Range.Text = "49+50="
$tbl.Cell(9, 2).Range.Text = "69+22="
$tbl.Cell(9, 3).Range.Text = "27-4="
$tbl.Cell(9, 4).Range.Text = "67-34="
$tbl.Cell(9, 5).Range.Text = "63-7="
$tbl.Cell(10, 1).Range.Text = "6+71="
$tbl.Cell(10, 2).Range.Text = "55-19="
$tbl.Cell(10, 3).Range.Text = "46+46="
$tbl.Cell(10, 4).Range.Text = "67-11="
$tbl.Cell(10, 5).Range.Text = "13+46="
$tbl.Cell(11, 1).Range.Text = "84-29="
$tbl.Cell(11, 2).Range.Text = "17+63="
$tbl.Cell(11, 3).Range.Text = "13+77="
$tbl.Cell(11, 4).Range.Text = "23+64="
$tbl.Cell(11, 5).Range.Text = "14+30="
$tbl.Cell(12, 1).Range.Text = "7-2="
$tbl.Cell(12, 2).Range.Text = "64-63="
$tbl.Cell(12, 3).Range.Text = "18+13="
$tbl.Cell(12, 4).Range.Text = "44-11="
$tbl.Cell(12, 5).Range.Text = "23+45="
$tbl.Cell(13, 1).Range.Text = "99-34="
$tbl.Cell(13, 2).Range.Text = "21-2="
$tbl.Cell(13, 3).Range.Text = "68+20="
$tbl.Cell(13, 4).Range.Text = "48+35="
$tbl.Cell(13, 5).Range.Text = "86-33="
$tbl.Cell(14, 1).Range.Text = "57+38="
$tbl.Cell(14, 2).Range.Text = "79-29="
$tbl.Cell(14, 3).Range.Text = "12+59="
$tbl.Cell(14, 4).Range.Text = "43+37="
$tbl.Cell(14, 5).Range.Text = "16+83="
$tbl.Cell(15, 1).Range.Text = "47+6="
$tbl.Cell(15, 2).Range.Text = "11+80="
$tbl.Cell(15, 3).Range.Text = "51+25="
$tbl.Cell(15, 4).Range.Text = "30+16="
$tbl.Cell(15, 5).Range.Text = "35-21="
$tbl.Cell(16, 1).Range.Text = "27-13="
$tbl.Cell(16, 2).Range.Text = "79-50="
$tbl.Cell(16, 3).Range.Text = "30+64="
$tbl.Cell(16, 4).Range.Text = "35+43="
$tbl.Cell(16, 5).Range.Text = "86+13="
$tbl.Cell(17, 1).Range.Text = "31-5="
$tbl.Cell(17, 2).Range.Text = "75+8="
$tbl.Cell(17, 3).Range.Text = "57+5="
$tbl.Cell(17, 4).Range.Text = "54-4="
$tbl.Cell(17, 5).Range.Text = "69-9="
$tbl.Cell(18, 1).Range.Text = "32+3="
$tbl.Cell(18, 2).Range.Text = "63-48="
$tbl.Cell(18, 3).Range.Text = "86-41="
$tbl.Cell(18, 4).Range.Text = "10+89="
$tbl.Cell(18, 5).Range.Text = "36+48="
$tbl.Cell(19, 1).Range.Text = "54-42="
$tbl.Cell(19, 2).Range.Text = "26+22="
$tbl.Cell(19, 3).Range.Text = "69-60="
$tbl.Cell(19, 4).Range.Text = "1+21="
$tbl.Cell(19, 5).Range.Text = "29+40="
$tbl.Cell(20, 1).Range.Text = "56-27="
$tbl.Cell(20, 2).Range.Text = "97-25="
$tbl.Cell(20, 3).Range.Text = "33+23="
$tbl.Cell(20, 4).Range.Text = "2+45="
$tbl.Cell(20, 5).Range.Text = "12+74="
